# Update the subtitle on the title slide: "SER 2020" -> "SER 2021"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Subtitle 2")
$tr = $shp.TextFrame.TextRange

$found = $tr.Find("SER 2020", 0, $false)
$found.Text = "SER 2021"
